# Update "想去人数" (want-to-go count) and "最低票价" (min ticket price)
# figures on the 展览 (Exhibition) and 全部类型 (All types) sheets to the
# freshly scraped values.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws = $wb.Worksheets.Item("展览")

$ws.Range("F2").Value = 448

$ws.Range("F3").Value = 1836
$ws.Range("G3").Value = 65

$ws.Range("F4").Value = 1448

$ws.Range("F5").Value = 140

$ws.Range("F6").Value = 1726

$ws.Range("F10").Value = 29

$ws.Range("F18").Value = 65

$ws.Range("F20").Value = 4524

$ws.Range("F22").Value = 808

$ws.Range("F23").Value = 97

$ws.Range("F24").Value = 2156

$ws.Range("F27").Value = 2023

# --- Sheet "全部类型" ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("全部类型")

$ws2.Range("F2").Value = 448

$ws2.Range("F3").Value = 1836
$ws2.Range("G3").Value = 65

$ws2.Range("F4").Value = 1448

$ws2.Range("F5").Value = 140

$ws2.Range("F6").Value = 1726

$ws2.Range("F10").Value = 29

$ws2.Range("F18").Value = 65

$ws2.Range("F20").Value = 4524

$ws2.Range("F24").Value = 808

$ws2.Range("F25").Value = 97

$ws2.Range("F26").Value = 2156

$ws2.Range("F29").Value = 2023
